$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3769309520721436
$ws.Range("E2").Value = 1347.237095632645
$ws.Range("F2").Value = 0.07990143729468036
$ws.Range("G2").Value = 0.05998581805288029
$ws.Range("H2").Value = 0.04830954899111276
$ws.Range("I2").Value = 0.04120592795978379
$ws.Range("J2").Value = 0.03753394998839247
$ws.Range("K2").Value = 0.03449438307065888
$ws.Range("L2").Value = 0.03198355730771082
$ws.Range("M2").Value = 0.03073307284517923
$ws.Range("N2").Value = 0.03001980329602668
$ws.Range("O2").Value = 0.02845043430031038
$ws.Range("P2").Value = 0.02825971302474251
$ws.Range("Q2").Value = 0.02763909695204155
$ws.Range("R2").Value = 0.02740726987440697
$ws.Range("S2").Value = 0.02710509825171897
$ws.Range("T2").Value = 0.02673387925186966
$ws.Range("U2").Value = 0.02656260428749077
$ws.Range("V2").Value = 0.02654760814120844
$ws.Range("W2").Value = 0.026415875730339
$ws.Range("X2").Value = 0.02630171111493943
$ws.Range("Y2").Value = 0.02626193168874551
$ws.Range("C3").Value = 0.3906493186950684
$ws.Range("E3").Value = 1377.17140992205
$ws.Range("F3").Value = 0.0812924136272716
$ws.Range("G3").Value = 0.05540204592206265
$ws.Range("H3").Value = 0.04843824428125919
$ws.Range("I3").Value = 0.04120704785610278
$ws.Range("J3").Value = 0.03538787365383474
$ws.Range("K3").Value = 0.03165305124293578
$ws.Range("L3").Value = 0.03038526235115781
$ws.Range("M3").Value = 0.02967211946265438
$ws.Range("N3").Value = 0.02892570186047997
$ws.Range("O3").Value = 0.02849763790545108
$ws.Range("P3").Value = 0.02844728464591295
$ws.Range("Q3").Value = 0.02808828639060565
$ws.Range("R3").Value = 0.02785642323902855
$ws.Range("S3").Value = 0.02752835968763482
$ws.Range("T3").Value = 0.02740285697124642
$ws.Range("U3").Value = 0.02714033014888262
$ws.Range("V3").Value = 0.02701369664808484
$ws.Range("W3").Value = 0.02696351456726834
$ws.Range("X3").Value = 0.02689671090913399
$ws.Range("Y3").Value = 0.02684544658717445
$ws.Range("C4").Value = 0.3749992847442627
$ws.Range("E4").Value = 1382.896236002034
$ws.Range("F4").Value = 0.08542942239929553
$ws.Range("G4").Value = 0.06420201572211295
$ws.Range("H4").Value = 0.05306851940032535
$ws.Range("I4").Value = 0.04313612718351405
$ws.Range("J4").Value = 0.03860979666464206
$ws.Range("K4").Value = 0.03546097305629935
$ws.Range("L4").Value = 0.03416837276371454
$ws.Range("M4").Value = 0.03238261866954559
$ws.Range("N4").Value = 0.0304731099026484
$ws.Range("O4").Value = 0.02998770616879475
$ws.Range("P4").Value = 0.02868310103030666
$ws.Range("Q4").Value = 0.02842361082705329
$ws.Range("R4").Value = 0.02827547796889425
$ws.Range("S4").Value = 0.02793726000656008
$ws.Range("T4").Value = 0.02756289970988156
$ws.Range("U4").Value = 0.02733678333022713
$ws.Range("V4").Value = 0.02720202100293223
$ws.Range("W4").Value = 0.02712972721020769
$ws.Range("X4").Value = 0.02697496256656123
$ws.Range("Y4").Value = 0.02695704163746655
$ws.Range("C5").Value = 0.3906257152557373
$ws.Range("E5").Value = 1393.067061104086
$ws.Range("F5").Value = 0.08271636652101659
$ws.Range("G5").Value = 0.06110280391988943
$ws.Range("H5").Value = 0.05115583417926205
$ws.Range("I5").Value = 0.04045364637849877
$ws.Range("J5").Value = 0.03836414959639275
$ws.Range("K5").Value = 0.03563850534883382
$ws.Range("L5").Value = 0.0336870162979081
$ws.Range("M5").Value = 0.03219659842881561
$ws.Range("N5").Value = 0.03071759968470434
$ws.Range("O5").Value = 0.02939913602080153
$ws.Range("P5").Value = 0.02923734832254393
$ws.Range("Q5").Value = 0.02850009791000295
$ws.Range("R5").Value = 0.02816624457430384
$ws.Range("S5").Value = 0.0278790850031571
$ws.Range("T5").Value = 0.02765985566322388
$ws.Range("U5").Value = 0.02746112650607684
$ws.Range("V5").Value = 0.02731408970998918
$ws.Range("W5").Value = 0.02721501677417059
$ws.Range("X5").Value = 0.02721501677417059
$ws.Range("Y5").Value = 0.02715530333536229
$ws.Range("C6").Value = 0.3593754768371582
$ws.Range("E6").Value = 1431.706987464699
$ws.Range("F6").Value = 0.07943900346020213
$ws.Range("G6").Value = 0.06195005432273458
$ws.Range("H6").Value = 0.05087821700165437
$ws.Range("I6").Value = 0.0427661908408841
$ws.Range("J6").Value = 0.03872280387127989
$ws.Range("K6").Value = 0.03770659284998746
$ws.Range("L6").Value = 0.03580434846694183
$ws.Range("M6").Value = 0.03335192709016408
$ws.Range("N6").Value = 0.03220677249516021
$ws.Range("O6").Value = 0.03115543324507835
$ws.Range("P6").Value = 0.03030140548679953
$ws.Range("Q6").Value = 0.02945433427796103
$ws.Range("R6").Value = 0.02907497991934088
$ws.Range("S6").Value = 0.02883561533040013
$ws.Range("T6").Value = 0.02856283294367909
$ws.Range("U6").Value = 0.0283839278519208
$ws.Range("V6").Value = 0.02829167092707771
$ws.Range("W6").Value = 0.02811511437820352
$ws.Range("X6").Value = 0.02802708863142551
$ws.Range("Y6").Value = 0.02790851827416566
$ws.Range("C7").Value = 0.3437490463256836
$ws.Range("E7").Value = 1369.960297288419
$ws.Range("F7").Value = 0.08281630118676042
$ws.Range("G7").Value = 0.06310120016123695
$ws.Range("H7").Value = 0.04877147537067913
$ws.Range("I7").Value = 0.0404463178672678
$ws.Range("J7").Value = 0.03649245753140582
$ws.Range("K7").Value = 0.03533139673662399
$ws.Range("L7").Value = 0.03329537517395878
$ws.Range("M7").Value = 0.03054936022941621
$ws.Range("N7").Value = 0.03054936022941621
$ws.Range("O7").Value = 0.03010387775616501
$ws.Range("P7").Value = 0.02900586672357407
$ws.Range("Q7").Value = 0.02808189832182107
$ws.Range("R7").Value = 0.02798620676274356
$ws.Range("S7").Value = 0.02753059442850906
$ws.Range("T7").Value = 0.02730587612062468
$ws.Range("U7").Value = 0.02704636813314557
$ws.Range("V7").Value = 0.02684934245587985
$ws.Range("W7").Value = 0.02680168157363776
$ws.Range("X7").Value = 0.02672766441269286
$ws.Range("Y7").Value = 0.02670487908944285
$ws.Range("C8").Value = 0.4689059257507324
$ws.Range("E8").Value = 1409.193635213785
$ws.Range("F8").Value = 0.08572550220393098
$ws.Range("G8").Value = 0.06729260996533787
$ws.Range("H8").Value = 0.05164888787419662
$ws.Range("I8").Value = 0.04461442785045112
$ws.Range("J8").Value = 0.03865988244582446
$ws.Range("K8").Value = 0.03588476873420547
$ws.Range("L8").Value = 0.0349581328547425
$ws.Range("M8").Value = 0.03278610943029901
$ws.Range("N8").Value = 0.03160274479877817
$ws.Range("O8").Value = 0.03055635371021232
$ws.Range("P8").Value = 0.02972506607393667
$ws.Range("Q8").Value = 0.0290174326721278
$ws.Range("R8").Value = 0.02866173125254625
$ws.Range("S8").Value = 0.02815457976925432
$ws.Range("T8").Value = 0.02801755087028648
$ws.Range("U8").Value = 0.02784030086663525
$ws.Range("V8").Value = 0.02781025993275956
$ws.Range("W8").Value = 0.02764438487261324
$ws.Range("X8").Value = 0.02752310726111899
$ws.Range("Y8").Value = 0.02746966150514201
$ws.Range("C9").Value = 0.3908798694610596
$ws.Range("E9").Value = 1387.295921304925
$ws.Range("F9").Value = 0.07784680508602365
$ws.Range("G9").Value = 0.06058814863156797
$ws.Range("H9").Value = 0.0450321306319223
$ws.Range("I9").Value = 0.04268102327037445
$ws.Range("J9").Value = 0.03903990008380216
$ws.Range("K9").Value = 0.03621739522642124
$ws.Range("L9").Value = 0.03343755680031001
$ws.Range("M9").Value = 0.03214033575496013
$ws.Range("N9").Value = 0.03050814135246261
$ws.Range("O9").Value = 0.02970045231989474
$ws.Range("P9").Value = 0.02913429069755801
$ws.Range("Q9").Value = 0.02868588405855576
$ws.Range("R9").Value = 0.0283956661321633
$ws.Range("S9").Value = 0.0279229800916227
$ws.Range("T9").Value = 0.0276609110058631
$ws.Range("U9").Value = 0.02756295481537359
$ws.Range("V9").Value = 0.02740226741505533
$ws.Range("W9").Value = 0.02729885741190023
$ws.Range("X9").Value = 0.02717102357719404
$ws.Range("Y9").Value = 0.02704280548352681
$ws.Range("C10").Value = 0.3906130790710449
$ws.Range("E10").Value = 1401.230195101383
$ws.Range("F10").Value = 0.07697621201859635
$ws.Range("G10").Value = 0.05859850658574604
$ws.Range("H10").Value = 0.04569061782630148
$ws.Range("I10").Value = 0.04275387805777491
$ws.Range("J10").Value = 0.03792906906406411
$ws.Range("K10").Value = 0.03418956519495158
$ws.Range("L10").Value = 0.03289130925226277
$ws.Range("M10").Value = 0.03068617172477781
$ws.Range("N10").Value = 0.03040480819569876
$ws.Range("O10").Value = 0.02960525000453131
$ws.Range("P10").Value = 0.02855234238838184
$ws.Range("Q10").Value = 0.02855234238838184
$ws.Range("R10").Value = 0.02817458382754375
$ws.Range("S10").Value = 0.0281038054175167
$ws.Range("T10").Value = 0.02784279183076844
$ws.Range("U10").Value = 0.02754603518044946
$ws.Range("V10").Value = 0.02744625102593914
$ws.Range("W10").Value = 0.02739874074850982
$ws.Range("X10").Value = 0.02731898906485605
$ws.Range("Y10").Value = 0.02731442875441293
$ws.Range("C11").Value = 0.3905973434448242
$ws.Range("E11").Value = 1429.323035413876
$ws.Range("F11").Value = 0.07756075925181811
$ws.Range("G11").Value = 0.06197917387669797
$ws.Range("H11").Value = 0.04586465237267182
$ws.Range("I11").Value = 0.04123770663513738
$ws.Range("J11").Value = 0.03837676293553268
$ws.Range("K11").Value = 0.03506266786025974
$ws.Range("L11").Value = 0.03243094082775808
$ws.Range("M11").Value = 0.0314899222932583
$ws.Range("N11").Value = 0.03066729472098168
$ws.Range("O11").Value = 0.02976738634013497
$ws.Range("P11").Value = 0.02941722932503209
$ws.Range("Q11").Value = 0.02925862352390542
$ws.Range("R11").Value = 0.02871195139957207
$ws.Range("S11").Value = 0.02845125183356308
$ws.Range("T11").Value = 0.02842757783225973
$ws.Range("U11").Value = 0.02824068325066394
$ws.Range("V11").Value = 0.02806691951606621
$ws.Range("W11").Value = 0.02796480686813367
$ws.Range("X11").Value = 0.02789520339224082
$ws.Range("Y11").Value = 0.02786204747395469
